$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = [double]"4.21848668766117"
$ws.Cells.Item(3, 2).Value = [double]"10.00000009951842"
$ws.Cells.Item(4, 2).Value = [double]"40"
$ws.Cells.Item(5, 2).Value = [double]"2.161876239173414"
$ws.Cells.Item(6, 2).Value = [double]"1.305407289332279"
$ws.Cells.Item(7, 2).Value = [double]"0.08376756518842321"
$ws.Cells.Item(8, 2).Value = [double]"4.621878021264386"
$ws.Cells.Item(9, 2).Value = [double]"1.191473075533368"
$ws.Cells.Item(10, 2).Value = [double]"1.115443443227179"
$ws.Cells.Item(11, 2).Value = [double]"-0.07602963230618887"
$ws.Cells.Item(12, 2).Value = [double]"0.03999999147900437"
$ws.Cells.Item(13, 2).Value = [double]"0.3999999900651737"
$ws.Cells.Item(14, 2).Value = [double]"0.02030793424814088"
$ws.Cells.Item(15, 2).Value = [double]"73691.38384536856"
$ws.Cells.Item(16, 2).Value = [double]"0.6420304121433735"
$ws.Cells.Item(17, 2).Value = [double]"0.1638258714174546"
$ws.Cells.Item(18, 2).Value = [double]"0.01982477821186242"
$ws.Cells.Item(19, 2).Value = [double]"0.00253246131388649"
$ws.Cells.Item(20, 2).Value = [double]"2.163573273794861e-16"
$ws.Cells.Item(21, 2).Value = [double]"-1.357404264707409e-14"
$ws.Cells.Item(22, 2).Value = [double]"3.849353761713425"
$ws.Cells.Item(23, 2).Value = [double]"0.9798169500729075"
$ws.Cells.Item(24, 2).Value = [double]"-33.37452705455155"
$ws.Cells.Item(25, 2).Value = [double]"0.4730196006306132"
$ws.Cells.Item(26, 2).Value = [double]"1.980612674413935"
$ws.Cells.Item(27, 2).Value = [double]"0.005847885527026722"
$ws.Cells.Item(28, 2).Value = [double]"0.9947996635701725"
$ws.Cells.Item(29, 2).Value = [double]"0.9798169500729075"
$ws.Cells.Item(30, 2).Value = [double]"0.000148175243828173"
$ws.Cells.Item(31, 2).Value = [double]"3.1321976209891e-10"
$ws.Cells.Item(32, 2).Value = [double]"0.08273918652096125"
$ws.Cells.Item(34, 2).Value = [double]"0.009307867455430474"
$ws.Cells.Item(35, 2).Value = [double]"0.001205047791969087"
$ws.Cells.Item(36, 2).Value = [double]"0.002935051688970926"
$ws.Cells.Item(37, 2).Value = [double]"0.004119632438805138"
$ws.Cells.Item(39, 2).Value = [double]"7.432276609990739e-06"
$ws.Cells.Item(41, 2).Value = [double]"83.76756518842321"
$ws.Cells.Item(42, 2).Value = [double]"0.1846758691915898"
$ws.Cells.Item(43, 2).Value = [double]"0.03243752254065602"
$ws.Cells.Item(44, 2).Value = [double]"0.00196265304297438"
$ws.Cells.Item(45, 2).Value = [double]"0.0002507136700747625"
$ws.Cells.Item(50, 2).Value = [double]"0.008111681105153184"
$ws.Cells.Item(52, 2).Value = [double]"-9.999207221378128e-09"
$ws.Cells.Item(53, 2).Value = [double]"0.006205004828772089"
$ws.Cells.Item(54, 2).Value = [double]"-0.04229463601223338"
$ws.Cells.Item(55, 2).Value = [double]"-0.2320985618833411"
$ws.Cells.Item(56, 2).Value = [double]"-7.167478029048528e-19"
$ws.Cells.Item(57, 2).Value = [double]"2.847475812096534e-18"
$ws.Cells.Item(58, 2).Value = [double]"1.072731241935686"
$ws.Cells.Item(59, 2).Value = [double]"-2.813834067537204e-17"
$ws.Cells.Item(60, 2).Value = [double]"0.2320985618833411"
$ws.Cells.Item(61, 2).Value = [double]"-7.167478029048528e-19"
$ws.Cells.Item(62, 2).Value = [double]"2.731485827453528e-13"
$ws.Cells.Item(63, 2).Value = [double]"2.046283823921632e-17"
$ws.Cells.Item(64, 2).Value = [double]"0.6007437797594768"
$ws.Cells.Item(65, 2).Value = [double]"-1.575784546274665e-17"
$ws.Cells.Item(66, 2).Value = [double]"0.129978285232879"
$ws.Cells.Item(67, 2).Value = [double]"-6.251858262565062e-19"
$ws.Cells.Item(68, 2).Value = [double]"5.994740924970615e-13"
$ws.Cells.Item(69, 2).Value = [double]"1.784878360322779e-17"
$ws.Cells.Item(70, 2).Value = [double]"0.1308007000671443"
$ws.Cells.Item(71, 2).Value = [double]"0.1012978618161968"
$ws.Cells.Item(72, 2).Value = [double]"2.445301999028072"
$ws.Cells.Item(73, 2).Value = [double]"1.191438890454768"
$ws.Cells.Item(74, 2).Value = [double]"3.956962828013333e-22"
$ws.Cells.Item(75, 2).Value = [double]"1.119257107314558e-19"
$ws.Cells.Item(76, 2).Value = [double]"-0.09781205912471523"
$ws.Cells.Item(77, 2).Value = [double]"1.928555967048174e-20"
$ws.Cells.Item(78, 2).Value = [double]"0.09294592024695689"
$ws.Cells.Item(79, 2).Value = [double]"-4.511424074042302e-05"
$ws.Cells.Item(80, 2).Value = [double]"-5.465140748049398e-06"
$ws.Cells.Item(81, 2).Value = [double]"0.017556734985501"
$ws.Cells.Item(82, 2).Value = [double]"-0.07999999002070281"
$ws.Cells.Item(83, 2).Value = [double]"-1.406411133743902e-05"
$ws.Cells.Item(84, 2).Value = [double]"0.02999999000332941"
$ws.Cells.Item(85, 2).Value = [double]"-1.014327321690334"
$ws.Cells.Item(86, 2).Value = [double]"-0.04153901388748159"
$ws.Cells.Item(87, 2).Value = [double]"-0.004854107401058848"
$ws.Cells.Item(88, 2).Value = [double]"-0.1804923826267631"
$ws.Cells.Item(89, 2).Value = [double]"-0.4258050299700366"
$ws.Cells.Item(90, 2).Value = [double]"-0.01122823203308724"
$ws.Cells.Item(91, 2).Value = [double]"-0.06187147829518037"
$ws.Cells.Item(92, 2).Value = [double]"-1.435964086494934"
$ws.Cells.Item(93, 2).Value = [double]"-0.2118806223445346"
$ws.Cells.Item(94, 2).Value = [double]"-2.551438248844172e-17"
$ws.Cells.Item(95, 2).Value = [double]"-2.232967113957864e-19"
$ws.Cells.Item(96, 2).Value = [double]"-2.233450486397225"
$ws.Cells.Item(97, 2).Value = [double]"2.901399978197982e-17"
$ws.Cells.Item(98, 2).Value = [double]"-0.002120961907992047"
$ws.Cells.Item(99, 2).Value = [double]"-0.0003402299471877423"
$ws.Cells.Item(100, 2).Value = [double]"0.1140151903228311"
$ws.Cells.Item(101, 2).Value = [double]"0.08924667301264325"
$ws.Cells.Item(102, 2).Value = [double]"-0.001333963242440948"
$ws.Cells.Item(103, 2).Value = [double]"-0.080176841064349"
